# Update row 2 (Cxcl5-Cxcr2, FAPs -> ECs) with new TPM-derived values,
# then remove row 3 (the FAPs/MuSCs -> ECs row) entirely so only one
# data row remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 3 (shifts nothing below it up, but removes the row
# and its data completely, matching the diff which drops the row node).
$ws.Rows.Item(3).Delete()

# Row 2 values that changed (columns A-H stay the same).
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 0.1905664181213333
$ws.Range("R2").Value = 1.715097763092
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
